# Refresh the cryptos price/volume snapshot (cells are plain text, not
# numbers, so any value that Excel would otherwise auto-parse as a number
# is entered with a leading apostrophe to force text, then the cell style
# is reset to "Normal" so no stray number-format style gets attached).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.130.35"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "1.855.36"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'233.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.61%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.4693"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.79%  "
$ws.Range("D8").Value = "'0.2810"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("D9").Value = "'0.06530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.53%  "
$ws.Range("D10").Value = "'19.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.57%  "
$ws.Range("D11").Value = "'0.07789"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.23%  "
$ws.Range("D12").Value = "'96.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.18%  "
$ws.Range("D13").Value = "1.857.28"
$ws.Range("E13").Value = "  -2.82%  "
$ws.Range("D14").Value = "'5.081"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "'0.6678"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "'282.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.61%  "
$ws.Range("D17").Value = "30.160.62"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'1.000"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "'5.456"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.21%  "
$ws.Range("D20").Value = "'12.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("D21").Value = "2.105.72"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("B23").Value = "ShibaInu"
$ws.Range("C23").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D23").Value = "'0.000007206"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").Value = "'6.128"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.83%  "
$ws.Range("D25").Value = "'167.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'9.287"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.73%  "
$ws.Range("D27").Value = "'18.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("D28").Value = "'1.914"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.14%  "
$ws.Range("E29").Value = "  -2.89%  "
$ws.Range("D30").Value = "'0.09612"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.74%  "
$ws.Range("D31").Value = "'4.400"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("D32").Value = "'1.468"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").Value = "'4.079"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("D34").Value = "'0.04661"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.57%  "
$ws.Range("D35").Value = "'0.6940"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.48%  "
$ws.Range("D36").Value = "'1.088"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").Value = "'1.000"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.30%  "
$ws.Range("D38").Value = "'2.705"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").Value = "'0.01847"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.10%  "
$ws.Range("D40").Value = "'6.271"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.96%  "
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("D42").Value = "'71.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.13%  "
$ws.Range("D43").Value = "'0.8580"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "
$ws.Range("D44").Value = "'1.943"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.98%  "
$ws.Range("B45").Value = "PaxDollar"
$ws.Range("C45").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").Value = "'103.93"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.67%  "
$ws.Range("D47").Value = "'0.4147"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.71%  "
$ws.Range("D48").Value = "1.015.16"
$ws.Range("E48").Value = "  +6.46%  "
$ws.Range("D49").Value = "'7.183"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("D50").Value = "'8.913"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").Value = "'33.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.90%  "
